$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.615.49'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '2.324.30'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.68'
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.02'
$ws.Range("E6").Value = '  -1.90%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("E8").Value = '  -0.95%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.101'
$ws.Range("E9").Value = '  -2.28%  '
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.23'
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.338'
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.53'
$ws.Range("E13").Value = '  -1.48%  '
$ws.Range("D14").Value = '2.733.81'
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("D15").Value = '56.570.50'
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000133'
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").Value = '2.312.37'
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.40'
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '329.99'
$ws.Range("E19").Value = '  +2.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.16'
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.78'
$ws.Range("E21").Value = '  +3.31%  '
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.14'
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.165'
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("E25").Value = '  +7.06%  '
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.31'
$ws.Range("E27").Value = '  +1.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.30'
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.69'
$ws.Range("E29").Value = '  -1.34%  '
$ws.Range("D30").Value = '0.0₃0721'
$ws.Range("E30").Value = '  -2.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.16'
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.32'
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("E34").Value = '  +0.47%  '
$ws.Range("E35").Value = '  -0.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.94'
$ws.Range("E36").Value = '  -2.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.884'
$ws.Range("E37").Value = '  -4.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.58'
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.58'
$ws.Range("E39").Value = '  +1.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '149.66'
$ws.Range("E40").Value = '  +7.91%  '
$ws.Range("E41").Value = '  -2.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.58'
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '280.52'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.07'
$ws.Range("E44").Value = '  -4.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0930'
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0497'
$ws.Range("E46").Value = '  -1.85%  '
$ws.Range("E47").Value = '  -1.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.42'
$ws.Range("E48").Value = '  +3.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0215'
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.17'
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("E51").Value = '  +0.74%  '
